$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 1186.8572
$ws.Range("I28").Value = 501.41666
$ws.Range("K28").Value = 501.41666
$ws.Range("M28").Value = -16.41665999999998
$ws.Range("H106").Value = 1954
$ws.Range("I106").Value = 1382.48
$ws.Range("J106").Value = 3995.1428
$ws.Range("K106").Value = 1382.48
$ws.Range("L106").Value = 3995.1428
$ws.Range("M106").Value = -751.48
$ws.Range("N106").Value = -5257.1428
$ws.Range("H107").Value = 838.3871
$ws.Range("I107").Value = 810.8
$ws.Range("J107").Value = 953.3333
$ws.Range("K107").Value = 810.8
$ws.Range("L107").Value = 953.3333
$ws.Range("M107").Value = 1109.2
$ws.Range("N107").Value = -4793.3333
$ws.Range("H132").Value = 11260.239
$ws.Range("I132").Value = 11514.231
$ws.Range("K132").Value = 34542.693
$ws.Range("M132").Value = -32012.693

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1608.881
$ws.Range("I2").Value = 1297.1562
$ws.Range("K2").Value = 1297.1562
$ws.Range("M2").Value = -1184.1562
$ws.Range("H45").Value = 3575.2144
$ws.Range("I45").Value = 2217.8
$ws.Range("K45").Value = 2217.8
$ws.Range("M45").Value = -1840.8
$ws.Range("H116").Value = 1608.881
$ws.Range("I116").Value = 1297.1562
$ws.Range("K116").Value = 1297.1562
$ws.Range("M116").Value = 996.8438000000001
$ws.Range("H132").Value = 1855.8837
$ws.Range("I132").Value = 1234.3704
$ws.Range("K132").Value = 3703.1112
$ws.Range("M132").Value = -1173.1112

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1608.881
$ws.Range("I3").Value = 1297.1562
$ws.Range("K3").Value = 1297.1562
$ws.Range("M3").Value = -1183.1562
$ws.Range("H57").Value = 100000
$ws.Range("J57").Value = 100000
$ws.Range("L57").Value = 100000
$ws.Range("N57").Value = -101440
$ws.Range("H99").Value = 3889.0715
$ws.Range("I99").Value = 1625.4
$ws.Range("K99").Value = 1625.4
$ws.Range("M99").Value = -127.4000000000001
$ws.Range("H134").Value = 2716.75
$ws.Range("I134").Value = 2538.7778
$ws.Range("K134").Value = 7616.3334
$ws.Range("M134").Value = -5081.3334
$ws.Range("H136").Value = 100000
$ws.Range("J136").Value = 100000
$ws.Range("L136").Value = 100000
$ws.Range("N136").Value = -110200
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()
$ws.Range("H140").Value = 65663.664
$ws.Range("J140").Value = 93539.73
$ws.Range("L140").Value = 93539.73
$ws.Range("N140").Value = -103899.73

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1324.5
$ws.Range("I16").Value = 1209.5
$ws.Range("K16").Value = 1209.5
$ws.Range("M16").Value = -922.5
$ws.Range("H31").Value = 2633927.8
$ws.Range("J31").Value = 2986.8
$ws.Range("L31").Value = 2986.8
$ws.Range("N31").Value = -3576.8
$ws.Range("H34").Value = 2633927.8
$ws.Range("J34").Value = 2986.8
$ws.Range("L34").Value = 2986.8
$ws.Range("N34").Value = -3390.8
$ws.Range("H86").Value = 39139.75
$ws.Range("I86").Value = 60488.8
$ws.Range("K86").Value = 60488.8
$ws.Range("M86").Value = -59365.8
$ws.Range("H89").Value = 39139.75
$ws.Range("I89").Value = 60488.8
$ws.Range("K89").Value = 302444
$ws.Range("M89").Value = -296828
$ws.Range("H107").Value = 808.88
$ws.Range("I107").Value = 672
$ws.Range("J107").Value = 1160.8572
$ws.Range("K107").Value = 672
$ws.Range("L107").Value = 1160.8572
$ws.Range("M107").Value = 1248
$ws.Range("N107").Value = -5000.8572
$ws.Range("H113").Value = 1324.5
$ws.Range("I113").Value = 1209.5
$ws.Range("K113").Value = 1209.5
$ws.Range("M113").Value = 960.5
$ws.Range("H132").Value = 21432.666
$ws.Range("I132").Value = 24249.424
$ws.Range("K132").Value = 72748.272
$ws.Range("M132").Value = -70218.272
$ws.Range("H134").Value = 2558
$ws.Range("I134").Value = 2271.1333
$ws.Range("K134").Value = 6813.3999
$ws.Range("M134").Value = -4278.3999

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 771706.6
$ws.Range("J2").Value = 248.2
$ws.Range("L2").Value = 1489.2
$ws.Range("N2").Value = -1715.2
$ws.Range("H7").Value = 4273595
$ws.Range("I7").Value = 112.71429
$ws.Range("K7").Value = 338.14287
$ws.Range("M7").Value = -226.14287
$ws.Range("H26").Value = 61.142857
$ws.Range("I26").Value = 44.8
$ws.Range("J26").Value = 102
$ws.Range("K26").Value = 134.4
$ws.Range("L26").Value = 306
$ws.Range("M26").Value = 153.6
$ws.Range("N26").Value = -882
$ws.Range("H33").Value = 119.1
$ws.Range("J33").Value = 50.75
$ws.Range("L33").Value = 304.5
$ws.Range("N33").Value = -870.5
$ws.Range("H92").Value = 1312.8667
$ws.Range("I92").Value = 2616.1667
$ws.Range("K92").Value = 7848.500100000001
$ws.Range("M92").Value = -6600.500100000001
$ws.Range("H117").Value = 1172.1111
$ws.Range("I117").Value = 1341.6666
$ws.Range("J117").Value = 833
$ws.Range("K117").Value = 4024.9998
$ws.Range("L117").Value = 2499
$ws.Range("M117").Value = -582.9998000000001
$ws.Range("N117").Value = -9383
$ws.Range("H121").Value = 153267.28
$ws.Range("I121").Value = 355043.34
$ws.Range("J121").Value = 1935.25
$ws.Range("K121").Value = 1065130.02
$ws.Range("L121").Value = 5805.75
$ws.Range("M121").Value = -1063820.02
$ws.Range("N121").Value = -8425.75
$ws.Range("H131").Value = 475424.12
$ws.Range("I131").Value = 2126061
$ws.Range("J131").Value = 3813.5715
$ws.Range("K131").Value = 6378183
$ws.Range("L131").Value = 11440.7145
$ws.Range("M131").Value = -6373143
$ws.Range("N131").Value = -21520.7145
$ws.Range("H137").Value = 3389.6667
$ws.Range("I137").Value = 3830.5
$ws.Range("J137").Value = 2508
$ws.Range("K137").Value = 11491.5
$ws.Range("L137").Value = 7524
$ws.Range("M137").Value = -6391.5
$ws.Range("N137").Value = -17724

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 15461.286
$ws.Range("I102").Value = 17436.7
$ws.Range("K102").Value = 17436.7
$ws.Range("M102").Value = -15814.7

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 8351.462
$ws.Range("I82").Value = 7040.68
$ws.Range("J82").Value = 10692.143
$ws.Range("K82").Value = 7040.68
$ws.Range("L82").Value = 10692.143
$ws.Range("M82").Value = -6679.68
$ws.Range("N82").Value = -11414.143
$ws.Range("H85").Value = 8351.462
$ws.Range("I85").Value = 7040.68
$ws.Range("J85").Value = 10692.143
$ws.Range("K85").Value = 7040.68
$ws.Range("L85").Value = 10692.143
$ws.Range("M85").Value = -5792.68
$ws.Range("N85").Value = -13188.143
$ws.Range("H104").Value = 17737
$ws.Range("J104").Value = 17737
$ws.Range("L104").Value = 17737
$ws.Range("N104").Value = -24725
$ws.Range("H132").Value = 6147.8823
$ws.Range("I132").Value = 6080.933
$ws.Range("K132").Value = 18242.799
$ws.Range("M132").Value = -15712.799
$ws.Range("H136").Value = 4861.1665
$ws.Range("I136").Value = 3833.4
$ws.Range("K136").Value = 11500.2
$ws.Range("M136").Value = -8950.200000000001

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 7528
$ws.Range("I62").Value = 4249.5
$ws.Range("J62").Value = 8464.714
$ws.Range("K62").Value = 4249.5
$ws.Range("L62").Value = 8464.714
$ws.Range("M62").Value = -3625.5
$ws.Range("N62").Value = -9712.714
$ws.Range("H65").Value = 7528
$ws.Range("I65").Value = 4249.5
$ws.Range("J65").Value = 8464.714
$ws.Range("K65").Value = 21247.5
$ws.Range("L65").Value = 42323.57
$ws.Range("M65").Value = -18127.5
$ws.Range("N65").Value = -48563.57
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()
$ws.Range("H125").Value = 79715
$ws.Range("J125").Value = 79715
$ws.Range("L125").Value = 79715
$ws.Range("N125").Value = -89555
$ws.Range("H126").Value = 457863.2
$ws.Range("I126").Value = 3593.75
$ws.Range("K126").Value = 10781.25
$ws.Range("M126").Value = -8311.25
$ws.Range("H132").Value = 20258.143
$ws.Range("I132").Value = 23732.943
$ws.Range("K132").Value = 71198.829
$ws.Range("M132").Value = -68668.829
